$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1960
$ws.Range("J17").Value = 1970
$ws.Range("L17").Value = 5910
$ws.Range("N17").Value = -6246
$ws.Range("H29").Value = 2166.5
$ws.Range("I29").Value = 999
$ws.Range("J29").Value = 2400
$ws.Range("K29").Value = 2997
$ws.Range("L29").Value = 7200
$ws.Range("M29").Value = -2716
$ws.Range("N29").Value = -7762
$ws.Range("H38").Value = 887.2
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 6000
$ws.Range("N38").Value = -6744
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7500
$ws.Range("K43").Value = 7500
$ws.Range("M43").Value = -7431
$ws.Range("H70").Value = 3833
$ws.Range("J70").Value = 3833
$ws.Range("L70").Value = 11499
$ws.Range("N70").Value = -12039
$ws.Range("H73").Value = 3833
$ws.Range("J73").Value = 3833
$ws.Range("L73").Value = 11499
$ws.Range("N73").Value = -13371
$ws.Range("H98").Value = 854.625
$ws.Range("I98").Value = 854.625
$ws.Range("K98").Value = 854.625
$ws.Range("M98").Value = 643.375
$ws.Range("H100").Value = 4878.8
$ws.Range("J100").Value = 4799.5
$ws.Range("L100").Value = 4799.5
$ws.Range("N100").Value = -5881.5
$ws.Range("H115").Value = 486
$ws.Range("I115").Value = 487.6
$ws.Range("J115").Value = 478
$ws.Range("K115").Value = 1462.8
$ws.Range("L115").Value = 1434
$ws.Range("M115").Value = 104.1999999999998
$ws.Range("N115").Value = -4568
$ws.Range("H122").Value = 854.625
$ws.Range("I122").Value = 854.625
$ws.Range("K122").Value = 2563.875
$ws.Range("M122").Value = -113.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4131.5713
$ws.Range("I32").Value = 2375.4707
$ws.Range("K32").Value = 2375.4707
$ws.Range("M32").Value = -2088.4707

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 535.3333
$ws.Range("I5").Value = 342.4
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 342.4
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -229.4
$ws.Range("N5").Value = -1726
$ws.Range("H94").Value = 13113.333
$ws.Range("I94").Value = 2251.25
$ws.Range("K94").Value = 2251.25
$ws.Range("M94").Value = -1800.25
$ws.Range("H99").Value = 5719.6
$ws.Range("I99").Value = 5719.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5719.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4221.6
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 4448.3335
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7034.4116
$ws.Range("I31").Value = 6122.5386
$ws.Range("J31").Value = 9998
$ws.Range("K31").Value = 6122.5386
$ws.Range("L31").Value = 9998
$ws.Range("M31").Value = -5827.5386
$ws.Range("N31").Value = -10588
$ws.Range("H34").Value = 7034.4116
$ws.Range("I34").Value = 6122.5386
$ws.Range("J34").Value = 9998
$ws.Range("K34").Value = 6122.5386
$ws.Range("L34").Value = 9998
$ws.Range("M34").Value = -5920.5386
$ws.Range("N34").Value = -10402
$ws.Range("H94").Value = 1451.5
$ws.Range("I94").Value = 1268.6666
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1268.6666
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -817.6666
$ws.Range("N94").Value = -2902
$ws.Range("H99").Value = 1912
$ws.Range("I99").Value = 1912
$ws.Range("K99").Value = 1912
$ws.Range("M99").Value = -414
$ws.Range("H105").Value = 948.3333
$ws.Range("I105").Value = 948.3333
$ws.Range("K105").Value = 948.3333
$ws.Range("M105").Value = 798.6667
$ws.Range("H122").Value = 1023.8333
$ws.Range("I122").Value = 1048.7646
$ws.Range("K122").Value = 3146.2938
$ws.Range("M122").Value = -696.2937999999999
$ws.Range("H126").Value = 1912
$ws.Range("I126").Value = 1912
$ws.Range("K126").Value = 5736
$ws.Range("M126").Value = -3266
$ws.Range("H132").Value = 5902.4
$ws.Range("I132").Value = 5753
$ws.Range("K132").Value = 17259
$ws.Range("M132").Value = -14729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3117.25
$ws.Range("J5").Value = 2962.5
$ws.Range("L5").Value = 8887.5
$ws.Range("N5").Value = -9111.5
$ws.Range("H14").Value = 735.6
$ws.Range("I14").Value = 735.6
$ws.Range("K14").Value = 2206.8
$ws.Range("M14").Value = -2033.8
$ws.Range("H23").Value = 1957
$ws.Range("I23").Value = 1995
$ws.Range("J23").Value = 1900
$ws.Range("K23").Value = 5985
$ws.Range("L23").Value = 5700
$ws.Range("M23").Value = -5750
$ws.Range("N23").Value = -6170
$ws.Range("H80").Value = 5687.5
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 6142.857
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 18428.571
$ws.Range("M80").Value = -6564
$ws.Range("N80").Value = -20300.571
$ws.Range("H83").Value = 5687.5
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 6142.857
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 55285.713
$ws.Range("M83").Value = -17820
$ws.Range("N83").Value = -64645.713
$ws.Range("H103").Value = 2999.3333
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 2999.2
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 8997.599999999999
$ws.Range("M103").Value = -8121
$ws.Range("N103").Value = -10755.6
$ws.Range("H135").Value = 3117.25
$ws.Range("J135").Value = 2962.5
$ws.Range("L135").Value = 26662.5
$ws.Range("N135").Value = -31732.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4250
$ws.Range("I29").Value = 4000
$ws.Range("K29").Value = 4000
$ws.Range("M29").Value = -3710
$ws.Range("H132").Value = 4790.8887
$ws.Range("I132").Value = 4790.8887
$ws.Range("K132").Value = 14372.6661
$ws.Range("M132").Value = -11842.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4708.6665
$ws.Range("I7").Value = 3778.2222
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 3778.2222
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -3666.2222
$ws.Range("N7").Value = -7724
$ws.Range("H9").Value = 939
$ws.Range("I9").Value = 909
$ws.Range("J9").Value = 999
$ws.Range("K9").Value = 909
$ws.Range("L9").Value = 999
$ws.Range("M9").Value = -685
$ws.Range("N9").Value = -1447
$ws.Range("H40").Value = 6955.5557
$ws.Range("I40").Value = 6800
$ws.Range("K40").Value = 6800
$ws.Range("M40").Value = -6664
$ws.Range("H126").Value = 4708.6665
$ws.Range("I126").Value = 3778.2222
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 11334.6666
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -8864.6666
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 23054
$ws.Range("I132").Value = 24935.75
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 74807.25
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -72277.25
$ws.Range("N132").Value = -29060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1889.7646
$ws.Range("I81").Value = 1385.6666
$ws.Range("J81").Value = 3099.6
$ws.Range("K81").Value = 2771.3332
$ws.Range("L81").Value = 6199.2
$ws.Range("M81").Value = -1710.3332
$ws.Range("N81").Value = -8321.200000000001
$ws.Range("H84").Value = 1889.7646
$ws.Range("I84").Value = 1385.6666
$ws.Range("J84").Value = 3099.6
$ws.Range("K84").Value = 13856.666
$ws.Range("L84").Value = 30996
$ws.Range("M84").Value = -8552.666000000001
$ws.Range("N84").Value = -41604

Write-Output "Applied all Kraken_Profits updates"